$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/25/2023  Through  12/31/2023"

# --- Data table updates (rows 15-30) ---
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "0"
$ws.Range("C14").Copy()
$c.PasteSpecial(-4122)

$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("M15").Value = -12.5
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 21
$ws.Range("H16").Value = 133.333333333333
$ws.Range("I16").Value = 220
$ws.Range("J16").Value = 233
$ws.Range("K16").Value = -5.57939914163
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = 53.846153846153
$ws.Range("N16").Value = -84.034833091436
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -8.333333333333
$ws.Range("I17").Value = 184
$ws.Range("J17").Value = 177
$ws.Range("K17").Value = 3.954802259887
$ws.Range("L17").Value = 17.948717948717
$ws.Range("M17").Value = 91.666666666666
$ws.Range("N17").Value = -37.627118644067
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 36
$ws.Range("H18").Value = -80.555555555555
$ws.Range("I18").Value = 251
$ws.Range("J18").Value = 287
$ws.Range("K18").Value = -12.543554006968
$ws.Range("L18").Value = -2.334630350194
$ws.Range("M18").Value = 7.725321888412
$ws.Range("N18").Value = -91.685988737992
$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 65
$ws.Range("F19").Value = 152
$ws.Range("G19").Value = 115
$ws.Range("H19").Value = 32.173913043478
$ws.Range("I19").Value = 1727
$ws.Range("J19").Value = 1737
$ws.Range("K19").Value = -0.575705238917
$ws.Range("L19").Value = 42.257001647446
$ws.Range("M19").Value = 33.668730650154
$ws.Range("N19").Value = -54.456751054852
$c = $ws.Range("C20")
$c.Value = 2
$ws.Range("F22").Copy()
$c.PasteSpecial(-4122)

$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -69.230769230769
$ws.Range("I20").Value = 151
$ws.Range("J20").Value = 190
$ws.Range("K20").Value = -20.526315789473
$ws.Range("L20").Value = -10.119047619047
$ws.Range("M20").Value = 71.590909090909
$ws.Range("N20").Value = -95.632050911194
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 197
$ws.Range("G21").Value = 186
$ws.Range("H21").Value = 5.913978494623
$ws.Range("I21").Value = 2550
$ws.Range("J21").Value = 2638
$ws.Range("K21").Value = -3.335860500379
$ws.Range("L21").Value = 28.334172118772
$ws.Range("M21").Value = 36.436597110754
$ws.Range("N21").Value = -78.707414829659
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = "0"
$ws.Range("C14").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0"
$ws.Range("C14").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "***.*"
$ws.Range("C14").Copy()
$c.PasteSpecial(-4122)

$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 150
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 58.620689655172
$c = $ws.Range("D23")
$c.Value = 1
$ws.Range("F22").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("E23")
$c.Value = 0
$ws.Range("H15").Copy()
$c.PasteSpecial(-4122)

$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = -15
$ws.Range("L23").Value = 9.677419354838
$ws.Range("M23").Value = 41.666666666666
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 69
$ws.Range("E24").Value = -40.579710144927
$ws.Range("F24").Value = 186
$ws.Range("G24").Value = 270
$ws.Range("H24").Value = -31.111111111111
$ws.Range("I24").Value = 3070
$ws.Range("J24").Value = 4019
$ws.Range("K24").Value = -23.61283901468
$ws.Range("L24").Value = 22.8
$ws.Range("M24").Value = 82.520808561236
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -23.529411764705
$ws.Range("I25").Value = 355
$ws.Range("J25").Value = 376
$ws.Range("K25").Value = -5.585106382978
$ws.Range("L25").Value = 3.197674418604
$ws.Range("M25").Value = 3.498542274052
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = "0"
$ws.Range("C14").Copy()
$c.PasteSpecial(-4122)

$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 88
$ws.Range("J27").Value = 85
$ws.Range("K27").Value = 3.529411764705
$ws.Range("L27").Value = -1.123595505617
$ws.Range("F30").Value = 6
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = "0"
$ws.Range("C14").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("H30")
$c.NumberFormat = "@"
$c.Value = "***.*"
$ws.Range("C14").Copy()
$c.PasteSpecial(-4122)

$ws.Range("I30").Value = 31
$ws.Range("K30").Value = 10.714285714285
$ws.Range("L30").Value = 106.666666666667

$excel.CutCopyMode = 0
